$d = $word.ActiveDocument

# Remove the "In progress" paragraph entirely (its text run plus the
# paragraph mark itself), leaving the preceding "F2025" paragraph intact.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*In progress*") {
        $p.Range.Delete()
    }
}
